$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39/40 swap: Aptos <-> TheSandbox (name, link, price and volume all move)
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"

# Price (D) and Volume(1h) (E) refresh for every coin row (2-51).
# Price column is stored as text in the workbook; force text via NumberFormat
# "@" before assigning so Excel does not auto-coerce values like "1.012" into a
# number, then ClearFormats() to drop the temporary number-format style again so
# the cell style matches the original (unstyled) cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.430.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.81"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.60%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4801"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.56%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4048"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08177"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.007"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.31"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.908.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.047"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.215"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06830"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.010"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.435.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.668"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.80"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.192"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.128.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.654"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "155.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.098"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.27"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.012"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09604"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.613"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.555"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.371"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06421"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02282"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.181"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5929"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.69"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.52%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.929"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1845"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.284"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.403"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07464"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5564"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.935"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.04"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.90%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.427"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.77"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.99%  "
